$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2

$rng = $ws.Range("A2:D40")
$key = $ws.Range("D2:D40")
$rng.Sort($key, 1)
